# App Finalizado e Funcionando 100%
#
# Fixes the CPF on row 17 (was mistakenly stored as text) and appends two
# new "diária" entries (rows 18-19) for 2025-03-01 and 2025-03-02 at
# "Cantina Volpi - Lauro".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: CPF was stored as text ("94585623502"); make it a real number ---
$ws.Cells.Item(17, 6).Value = 94585623502

# --- Row 18: 2025-03-01, Cantina Volpi - Lauro, Edmilson ---
$ws.Cells.Item(18, 1).Value = "2025-03-01 10:35:00"
$ws.Cells.Item(18, 2).Value = "2025-03-01 15:35:00"
$ws.Cells.Item(18, 3).Value = "Cantina Volpi - Lauro"
$ws.Cells.Item(18, 4).Value = "Cantina Volpi - Lauro"
$ws.Cells.Item(18, 5).Value = "Edmilson"
$ws.Cells.Item(18, 6).Value = 1234532608
$ws.Cells.Item(18, 7).Value = 50
$ws.Cells.Item(18, 8).Value = 50
$ws.Cells.Item(18, 9).Value = "N"
$ws.Cells.Item(18, 10).Value = "N"

# --- Row 19: 2025-03-02, Cantina Volpi - Lauro, new deliverer Carlos ---
$ws.Cells.Item(19, 1).Value = "2025-03-02 10:37:00"
$ws.Cells.Item(19, 2).Value = "2025-03-02 16:37:00"
$ws.Cells.Item(19, 3).Value = "Cantina Volpi - Lauro"
$ws.Cells.Item(19, 4).Value = "Cantina Volpi - Lauro"
$ws.Cells.Item(19, 5).Value = "Carlos"

# Carlos' CPF is short/invalid, so the original author kept it as text;
# force text storage (temporarily) then restore the default style so the
# cell keeps its plain look while the underlying value stays a string.
$ws.Cells.Item(19, 6).NumberFormat = "@"
$ws.Cells.Item(19, 6).Value = "1586485914"
$ws.Cells.Item(19, 6).Style = "Normal"

$ws.Cells.Item(19, 7).Value = 50
$ws.Cells.Item(19, 8).Value = 50
$ws.Cells.Item(19, 9).Value = "N"
$ws.Cells.Item(19, 10).Value = "N"
